$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Invalid / Absent marked
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

# Row 4: Total Attendance Count / Real marked
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1

# Rows 5-18: Absent marked
foreach ($r in 5..18) {
    $ws.Range("H$r").Value = 1
}
